$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (it will be re-added at its new location).
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBm.Delete()
} catch {
}

# 2. Locate the sentence "mechanical deformation reveals" and insert "(DIAMOND) "
#    right before "reveals", turning it into "mechanical deformation (DIAMOND) reveals".
$range = $d.Content
$found = $range.Find.Execute("mechanical deformation reveals", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $prefixLen = "mechanical deformation ".Length
    $insertPos = $range.Start + $prefixLen

    $insertionPoint = $d.Range($insertPos, $insertPos)
    $insertionPoint.InsertBefore("(DIAMOND) ")

    # 3. Re-create the "_GoBack" bookmark right after the newly inserted "(DIAMOND"
    #    text (i.e. immediately before the closing ") ").
    $bmPos = $insertPos + "(DIAMOND".Length
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
